$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recomputed statistics) ---
$ws.Range("G9").Value = 0.0812173141466191
$ws.Range("G10").Value = 0.0812173141466191
$ws.Range("G13").Value = 0.6832
$ws.Range("G14").Value = 0.6832
$ws.Range("G27").Value = 0.08063728775707719
$ws.Range("G28").Value = 0.08063728775707719
$ws.Range("G31").Value = 0.658130909090909
$ws.Range("G32").Value = 0.658130909090909
$ws.Range("G44").Value = 0.0765119683335705
$ws.Range("G45").Value = 0.0765119683335705
$ws.Range("G48").Value = 0.582898181818182
$ws.Range("L48").Value = 0.6001
$ws.Range("M48").Value = 0.80372
$ws.Range("G49").Value = 0.582898181818182
$ws.Range("L49").Value = 0.6001
$ws.Range("M49").Value = 0.80372
$ws.Range("G55").Value = 1.51152290097887
$ws.Range("H55").Value = 2.74873283132386
$ws.Range("G62").Value = 0.0535547958097898
$ws.Range("G63").Value = 0.0535547958097898
$ws.Range("G66").Value = 0.492682456140351
$ws.Range("I66").Value = 0.86687
$ws.Range("L66").Value = 0.5658
$ws.Range("M66").Value = 0.70105
$ws.Range("N66").Value = 0.81699
$ws.Range("G67").Value = 0.492682456140351
$ws.Range("I67").Value = 0.86687
$ws.Range("L67").Value = 0.5658
$ws.Range("M67").Value = 0.70105
$ws.Range("N67").Value = 0.81699
$ws.Range("G73").Value = 1.60176238029036
$ws.Range("I73").Value = 2.94193
$ws.Range("N73").Value = 2.52632
$ws.Range("F80").Value = 0.00165
$ws.Range("G80").Value = 0.0451101588820846
$ws.Range("L80").Value = 0.0009700000000000001
$ws.Range("F81").Value = 0.00165
$ws.Range("G81").Value = 0.0451101588820846
$ws.Range("L81").Value = 0.0009700000000000001
$ws.Range("F84").Value = 0.3764
$ws.Range("G84").Value = 0.426050877192982
$ws.Range("I84").Value = 0.76893
$ws.Range("M84").Value = 0.66049
$ws.Range("N84").Value = 0.70122
$ws.Range("F85").Value = 0.3764
$ws.Range("G85").Value = 0.426050877192982
$ws.Range("I85").Value = 0.76893
$ws.Range("M85").Value = 0.66049
$ws.Range("N85").Value = 0.70122
$ws.Range("F91").Value = 1.64
$ws.Range("G91").Value = 1.68871017297431
$ws.Range("I91").Value = 3.07322
$ws.Range("L91").Value = 1.85
$ws.Range("M91").Value = 2.26356
$ws.Range("N91").Value = 2.73834
$ws.Range("G94").Value = 936.8077107275629
$ws.Range("G95").Value = 936.8077107275629
$ws.Range("G96").Value = 936.8077107275629
$ws.Range("G97").Value = 936.8077107275629
$ws.Range("F98").Value = 0.00152
$ws.Range("G98").Value = 0.0239977700839333
$ws.Range("L98").Value = 0.00102
$ws.Range("F99").Value = 0.00152
$ws.Range("G99").Value = 0.0239977700839333
$ws.Range("L99").Value = 0.00102
$ws.Range("F102").Value = 0.3277
$ws.Range("G102").Value = 0.367598333333333
$ws.Range("H102").Value = 0.8425
$ws.Range("I102").Value = 0.7051500000000001
$ws.Range("L102").Value = 0.33395
$ws.Range("M102").Value = 0.54312
$ws.Range("N102").Value = 0.6831199999999999
$ws.Range("F103").Value = 0.3277
$ws.Range("G103").Value = 0.367598333333333
$ws.Range("H103").Value = 0.8425
$ws.Range("I103").Value = 0.7051500000000001
$ws.Range("L103").Value = 0.33395
$ws.Range("M103").Value = 0.54312
$ws.Range("N103").Value = 0.6831199999999999
$ws.Range("C109").Value = "B"
$ws.Range("F109").Value = 1.57
$ws.Range("G109").Value = 1.53347828716342
$ws.Range("I109").Value = 3.06786
$ws.Range("L109").Value = 1.83
$ws.Range("M109").Value = 2.13968
$ws.Range("N109").Value = 2.7134
$ws.Range("G112").Value = 872.5743773942301
$ws.Range("G113").Value = 872.5743773942301
$ws.Range("G114").Value = 872.5743773942301
$ws.Range("G115").Value = 872.5743773942301
$ws.Range("F116").Value = 0.00152
$ws.Range("G116").Value = 0.0132653868254879
$ws.Range("L116").Value = 0.00108
$ws.Range("M116").Value = 0.00381
$ws.Range("F117").Value = 0.00152
$ws.Range("G117").Value = 0.0132653868254879
$ws.Range("L117").Value = 0.00108
$ws.Range("M117").Value = 0.00381
$ws.Range("G120").Value = 0.311195
$ws.Range("H120").Value = 0.8425
$ws.Range("I120").Value = 0.6677999999999999
$ws.Range("M120").Value = 0.39408
$ws.Range("N120").Value = 0.53088
$ws.Range("G121").Value = 0.311195
$ws.Range("H121").Value = 0.8425
$ws.Range("I121").Value = 0.6677999999999999
$ws.Range("M121").Value = 0.39408
$ws.Range("N121").Value = 0.53088
$ws.Range("F127").Value = 1.615
$ws.Range("G127").Value = 1.5546739445961
$ws.Range("I127").Value = 3.08214
$ws.Range("L127").Value = 1.85
$ws.Range("M127").Value = 2.328
$ws.Range("N127").Value = 2.8434
$ws.Range("G130").Value = 240.770144116563
$ws.Range("G131").Value = 240.770144116563
$ws.Range("G132").Value = 240.770144116563
$ws.Range("G133").Value = 240.770144116563
$ws.Range("F134").Value = 0.00193
$ws.Range("G134").Value = 0.0030347411019031
$ws.Range("I134").Value = 0.009730000000000001
$ws.Range("L134").Value = 0.00264
$ws.Range("M134").Value = 0.00614
$ws.Range("N134").Value = 0.008489999999999999
$ws.Range("F135").Value = 0.00193
$ws.Range("G135").Value = 0.0030347411019031
$ws.Range("I135").Value = 0.009730000000000001
$ws.Range("L135").Value = 0.00264
$ws.Range("M135").Value = 0.00614
$ws.Range("N135").Value = 0.008489999999999999
$ws.Range("G138").Value = 0.282616666666667
$ws.Range("G139").Value = 0.282616666666667
$ws.Range("G147").Value = 240.154867925907
$ws.Range("G148").Value = 240.154867925907
$ws.Range("G149").Value = 240.154867925907
$ws.Range("G150").Value = 240.154867925907
$ws.Range("F151").Value = 0.00303
$ws.Range("G151").Value = 0.0033571842234038
$ws.Range("I151").Value = 0.009169999999999999
$ws.Range("L151").Value = 0.00302
$ws.Range("M151").Value = 0.00618
$ws.Range("N151").Value = 0.0075
$ws.Range("F152").Value = 0.00303
$ws.Range("G152").Value = 0.0033571842234038
$ws.Range("I152").Value = 0.009169999999999999
$ws.Range("L152").Value = 0.00302
$ws.Range("M152").Value = 0.00618
$ws.Range("N152").Value = 0.0075
$ws.Range("G155").Value = 0.286101694915254
$ws.Range("G156").Value = 0.286101694915254
$ws.Range("G164").Value = 69.71418995980569
$ws.Range("G165").Value = 69.71418995980569
$ws.Range("G166").Value = 69.71418995980569
$ws.Range("G167").Value = 69.71418995980569
$ws.Range("F168").Value = 0.0032
$ws.Range("G168").Value = 0.0035545762501921
$ws.Range("I168").Value = 0.00912
$ws.Range("L168").Value = 0.00314
$ws.Range("M168").Value = 0.00616
$ws.Range("N168").Value = 0.00738
$ws.Range("F169").Value = 0.0032
$ws.Range("G169").Value = 0.0035545762501921
$ws.Range("I169").Value = 0.00912
$ws.Range("L169").Value = 0.00314
$ws.Range("M169").Value = 0.00616
$ws.Range("N169").Value = 0.00738
$ws.Range("G172").Value = 0.277254237288136
$ws.Range("G173").Value = 0.277254237288136
$ws.Range("G181").Value = 73.6760636974906
$ws.Range("G182").Value = 73.6760636974906
$ws.Range("G183").Value = 73.6760636974906
$ws.Range("G184").Value = 73.6760636974906
$ws.Range("F185").Value = 0.00333
$ws.Range("G185").Value = 0.0040345051959891
$ws.Range("I185").Value = 0.0098
$ws.Range("L185").Value = 0.00347
$ws.Range("M185").Value = 0.0066
$ws.Range("N185").Value = 0.009220000000000001
$ws.Range("F186").Value = 0.00333
$ws.Range("G186").Value = 0.0040345051959891
$ws.Range("I186").Value = 0.0098
$ws.Range("L186").Value = 0.00347
$ws.Range("M186").Value = 0.0066
$ws.Range("N186").Value = 0.009220000000000001
$ws.Range("G189").Value = 0.310237288135593
$ws.Range("N189").Value = 0.38951
$ws.Range("G190").Value = 0.310237288135593
$ws.Range("N190").Value = 0.38951

# --- Append new rows 195-211 (2019-2023 period results) ---
$ws.Range("A195").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B195").Value = "Chlorophyll A (83rd Percentile)"
$ws.Range("C195").Value = "B"
$ws.Range("D195").Value = "2019 - 2023"
$ws.Range("E195").Value = "Impact"
$ws.Range("F195").Value = 19
$ws.Range("G195").Value = 29.3025423728814
$ws.Range("H195").Value = 120
$ws.Range("I195").Value = 87.75
$ws.Range("L195").Value = 27
$ws.Range("M195").Value = 60
$ws.Range("N195").Value = 80
$ws.Range("O195").Value = 1828659.395
$ws.Range("P195").Value = 5628533.029
$ws.Range("Q195").Value = "Ruapehu District"
$ws.Range("R195").Value = "Whangaehu"
$ws.Range("S195").Value = "Upper Whangaehu"
$ws.Range("T195").Value = "Whau_1b"
$ws.Range("U195").Value = "mg chl-a /m2"
$ws.Range("A196").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B196").Value = "DRP (95th Percentile)"
$ws.Range("C196").Value = "D"
$ws.Range("D196").Value = "2019 - 2023"
$ws.Range("E196").Value = "Impact"
$ws.Range("F196").Value = 0.03
$ws.Range("G196").Value = 0.0333898305084746
$ws.Range("H196").Value = 0.081
$ws.Range("I196").Value = 0.06725
$ws.Range("L196").Value = 0.029
$ws.Range("M196").Value = 0.047
$ws.Range("N196").Value = 0.053
$ws.Range("O196").Value = 1828659.395
$ws.Range("P196").Value = 5628533.029
$ws.Range("Q196").Value = "Ruapehu District"
$ws.Range("R196").Value = "Whangaehu"
$ws.Range("S196").Value = "Upper Whangaehu"
$ws.Range("T196").Value = "Whau_1b"
$ws.Range("U196").Value = "mg/L"
$ws.Range("A197").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B197").Value = "DRP (Median)"
$ws.Range("C197").Value = "D"
$ws.Range("D197").Value = "2019 - 2023"
$ws.Range("E197").Value = "Impact"
$ws.Range("F197").Value = 0.03
$ws.Range("G197").Value = 0.0333898305084746
$ws.Range("H197").Value = 0.081
$ws.Range("I197").Value = 0.06725
$ws.Range("L197").Value = 0.029
$ws.Range("M197").Value = 0.047
$ws.Range("N197").Value = 0.053
$ws.Range("O197").Value = 1828659.395
$ws.Range("P197").Value = 5628533.029
$ws.Range("Q197").Value = "Ruapehu District"
$ws.Range("R197").Value = "Whangaehu"
$ws.Range("S197").Value = "Upper Whangaehu"
$ws.Range("T197").Value = "Whau_1b"
$ws.Range("U197").Value = "mg/L"
$ws.Range("A198").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B198").Value = "E coli (>260)"
$ws.Range("C198").Value = "A"
$ws.Range("D198").Value = "2019 - 2023"
$ws.Range("E198").Value = "Impact"
$ws.Range("F198").Value = 16
$ws.Range("G198").Value = 27.5913179347787
$ws.Range("H198").Value = 295
$ws.Range("I198").Value = 78.2
$ws.Range("J198").Value = 0
$ws.Range("K198").Value = 1.69491525423729
$ws.Range("L198").Value = 27
$ws.Range("M198").Value = 45.88
$ws.Range("N198").Value = 67.12
$ws.Range("O198").Value = 1828659.395
$ws.Range("P198").Value = 5628533.029
$ws.Range("Q198").Value = "Ruapehu District"
$ws.Range("R198").Value = "Whangaehu"
$ws.Range("S198").Value = "Upper Whangaehu"
$ws.Range("T198").Value = "Whau_1b"
$ws.Range("U198").Value = "% exceedances over 260/100 mL"
$ws.Range("A199").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B199").Value = "E coli (>540)"
$ws.Range("C199").Value = "A"
$ws.Range("D199").Value = "2019 - 2023"
$ws.Range("E199").Value = "Impact"
$ws.Range("F199").Value = 16
$ws.Range("G199").Value = 27.5913179347787
$ws.Range("H199").Value = 295
$ws.Range("I199").Value = 78.2
$ws.Range("J199").Value = 0
$ws.Range("K199").Value = 1.69491525423729
$ws.Range("L199").Value = 27
$ws.Range("M199").Value = 45.88
$ws.Range("N199").Value = 67.12
$ws.Range("O199").Value = 1828659.395
$ws.Range("P199").Value = 5628533.029
$ws.Range("Q199").Value = "Ruapehu District"
$ws.Range("R199").Value = "Whangaehu"
$ws.Range("S199").Value = "Upper Whangaehu"
$ws.Range("T199").Value = "Whau_1b"
$ws.Range("U199").Value = "% exceedances over 540/100 mL"
$ws.Range("A200").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B200").Value = "E coli (Median)"
$ws.Range("C200").Value = "A"
$ws.Range("D200").Value = "2019 - 2023"
$ws.Range("E200").Value = "Impact"
$ws.Range("F200").Value = 16
$ws.Range("G200").Value = 27.5913179347787
$ws.Range("H200").Value = 295
$ws.Range("I200").Value = 78.2
$ws.Range("J200").Value = 0
$ws.Range("K200").Value = 1.69491525423729
$ws.Range("L200").Value = 27
$ws.Range("M200").Value = 45.88
$ws.Range("N200").Value = 67.12
$ws.Range("O200").Value = 1828659.395
$ws.Range("P200").Value = 5628533.029
$ws.Range("Q200").Value = "Ruapehu District"
$ws.Range("R200").Value = "Whangaehu"
$ws.Range("S200").Value = "Upper Whangaehu"
$ws.Range("T200").Value = "Whau_1b"
$ws.Range("U200").Value = "E. coli/100 mL"
$ws.Range("A201").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B201").Value = "E coli (95th Percentile)"
$ws.Range("C201").Value = "A"
$ws.Range("D201").Value = "2019 - 2023"
$ws.Range("E201").Value = "Impact"
$ws.Range("F201").Value = 16
$ws.Range("G201").Value = 27.5913179347787
$ws.Range("H201").Value = 295
$ws.Range("I201").Value = 78.2
$ws.Range("J201").Value = 0
$ws.Range("K201").Value = 1.69491525423729
$ws.Range("L201").Value = 27
$ws.Range("M201").Value = 45.88
$ws.Range("N201").Value = 67.12
$ws.Range("O201").Value = 1828659.395
$ws.Range("P201").Value = 5628533.029
$ws.Range("Q201").Value = "Ruapehu District"
$ws.Range("R201").Value = "Whangaehu"
$ws.Range("S201").Value = "Upper Whangaehu"
$ws.Range("T201").Value = "Whau_1b"
$ws.Range("U201").Value = "E. coli/100 mL"
$ws.Range("A202").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B202").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C202").Value = "A"
$ws.Range("D202").Value = "2019 - 2023"
$ws.Range("E202").Value = "Impact"
$ws.Range("F202").Value = 0.00324
$ws.Range("G202").Value = 0.0041300206583638
$ws.Range("H202").Value = 0.0130788892630831
$ws.Range("I202").Value = 0.01164
$ws.Range("L202").Value = 0.0031
$ws.Range("M202").Value = 0.00667
$ws.Range("N202").Value = 0.0094
$ws.Range("O202").Value = 1828659.395
$ws.Range("P202").Value = 5628533.029
$ws.Range("Q202").Value = "Ruapehu District"
$ws.Range("R202").Value = "Whangaehu"
$ws.Range("S202").Value = "Upper Whangaehu"
$ws.Range("T202").Value = "Whau_1b"
$ws.Range("U202").Value = "mg NH4-N/L"
$ws.Range("A203").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B203").Value = "Ammoniacal-N (Median)"
$ws.Range("C203").Value = "A"
$ws.Range("D203").Value = "2019 - 2023"
$ws.Range("E203").Value = "Impact"
$ws.Range("F203").Value = 0.00324
$ws.Range("G203").Value = 0.0041300206583638
$ws.Range("H203").Value = 0.0130788892630831
$ws.Range("I203").Value = 0.01164
$ws.Range("L203").Value = 0.0031
$ws.Range("M203").Value = 0.00667
$ws.Range("N203").Value = 0.0094
$ws.Range("O203").Value = 1828659.395
$ws.Range("P203").Value = 5628533.029
$ws.Range("Q203").Value = "Ruapehu District"
$ws.Range("R203").Value = "Whangaehu"
$ws.Range("S203").Value = "Upper Whangaehu"
$ws.Range("T203").Value = "Whau_1b"
$ws.Range("U203").Value = "mg NH4-N/L"
$ws.Range("A204").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B204").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C204").Value = "A"
$ws.Range("D204").Value = "2019 - 2023"
$ws.Range("E204").Value = "Impact"
$ws.Range("F204").Value = 0.267
$ws.Range("G204").Value = 0.304423728813559
$ws.Range("H204").Value = 1.94
$ws.Range("I204").Value = 0.4354
$ws.Range("L204").Value = 0.3205
$ws.Range("M204").Value = 0.36129
$ws.Range("N204").Value = 0.40456
$ws.Range("O204").Value = 1828659.395
$ws.Range("P204").Value = 5628533.029
$ws.Range("Q204").Value = "Ruapehu District"
$ws.Range("R204").Value = "Whangaehu"
$ws.Range("S204").Value = "Upper Whangaehu"
$ws.Range("T204").Value = "Whau_1b"
$ws.Range("U204").Value = "mg NO3-N/L"
$ws.Range("A205").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B205").Value = "Nitrate-N (Median)"
$ws.Range("C205").Value = "A"
$ws.Range("D205").Value = "2019 - 2023"
$ws.Range("E205").Value = "Impact"
$ws.Range("F205").Value = 0.267
$ws.Range("G205").Value = 0.304423728813559
$ws.Range("H205").Value = 1.94
$ws.Range("I205").Value = 0.4354
$ws.Range("L205").Value = 0.3205
$ws.Range("M205").Value = 0.36129
$ws.Range("N205").Value = 0.40456
$ws.Range("O205").Value = 1828659.395
$ws.Range("P205").Value = 5628533.029
$ws.Range("Q205").Value = "Ruapehu District"
$ws.Range("R205").Value = "Whangaehu"
$ws.Range("S205").Value = "Upper Whangaehu"
$ws.Range("T205").Value = "Whau_1b"
$ws.Range("U205").Value = "mg NO3-N/L"
$ws.Range("A206").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B206").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D206").Value = "2019 - 2023"
$ws.Range("E206").Value = "Impact"
$ws.Range("F206").Value = 0.275
$ws.Range("G206").Value = 0.311762711864407
$ws.Range("H206").Value = 1.943
$ws.Range("I206").Value = 0.44305
$ws.Range("L206").Value = 0.3325
$ws.Range("M206").Value = 0.36923
$ws.Range("N206").Value = 0.41882
$ws.Range("O206").Value = 1828659.395
$ws.Range("P206").Value = 5628533.029
$ws.Range("Q206").Value = "Ruapehu District"
$ws.Range("R206").Value = "Whangaehu"
$ws.Range("S206").Value = "Upper Whangaehu"
$ws.Range("T206").Value = "Whau_1b"
$ws.Range("U206").Value = "g/m3"
$ws.Range("A207").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B207").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D207").Value = "2019 - 2023"
$ws.Range("E207").Value = "Impact"
$ws.Range("F207").Value = 0.275
$ws.Range("G207").Value = 0.311762711864407
$ws.Range("H207").Value = 1.943
$ws.Range("I207").Value = 0.44305
$ws.Range("L207").Value = 0.3325
$ws.Range("M207").Value = 0.36923
$ws.Range("N207").Value = 0.41882
$ws.Range("O207").Value = 1828659.395
$ws.Range("P207").Value = 5628533.029
$ws.Range("Q207").Value = "Ruapehu District"
$ws.Range("R207").Value = "Whangaehu"
$ws.Range("S207").Value = "Upper Whangaehu"
$ws.Range("T207").Value = "Whau_1b"
$ws.Range("U207").Value = "g/m3"
$ws.Range("A208").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B208").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D208").Value = "2019 - 2023"
$ws.Range("E208").Value = "Impact"
$ws.Range("F208").Value = 0.34
$ws.Range("G208").Value = 0.350677966101695
$ws.Range("H208").Value = 0.6899999999999999
$ws.Range("I208").Value = 0.501
$ws.Range("L208").Value = 0.375
$ws.Range("M208").Value = 0.4247
$ws.Range("N208").Value = 0.46
$ws.Range("O208").Value = 1828659.395
$ws.Range("P208").Value = 5628533.029
$ws.Range("Q208").Value = "Ruapehu District"
$ws.Range("R208").Value = "Whangaehu"
$ws.Range("S208").Value = "Upper Whangaehu"
$ws.Range("T208").Value = "Whau_1b"
$ws.Range("U208").Value = "g/m3"
$ws.Range("A209").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B209").Value = "Total Nitrogen (Median)"
$ws.Range("D209").Value = "2019 - 2023"
$ws.Range("E209").Value = "Impact"
$ws.Range("F209").Value = 0.34
$ws.Range("G209").Value = 0.350677966101695
$ws.Range("H209").Value = 0.6899999999999999
$ws.Range("I209").Value = 0.501
$ws.Range("L209").Value = 0.375
$ws.Range("M209").Value = 0.4247
$ws.Range("N209").Value = 0.46
$ws.Range("O209").Value = 1828659.395
$ws.Range("P209").Value = 5628533.029
$ws.Range("Q209").Value = "Ruapehu District"
$ws.Range("R209").Value = "Whangaehu"
$ws.Range("S209").Value = "Upper Whangaehu"
$ws.Range("T209").Value = "Whau_1b"
$ws.Range("U209").Value = "g/m3"
$ws.Range("A210").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B210").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D210").Value = "2019 - 2023"
$ws.Range("E210").Value = "Impact"
$ws.Range("F210").Value = 0.037
$ws.Range("G210").Value = 0.0489152542372881
$ws.Range("H210").Value = 0.317
$ws.Range("I210").Value = 0.09594999999999999
$ws.Range("L210").Value = 0.035
$ws.Range("M210").Value = 0.05782
$ws.Range("N210").Value = 0.07056
$ws.Range("O210").Value = 1828659.395
$ws.Range("P210").Value = 5628533.029
$ws.Range("Q210").Value = "Ruapehu District"
$ws.Range("R210").Value = "Whangaehu"
$ws.Range("S210").Value = "Upper Whangaehu"
$ws.Range("T210").Value = "Whau_1b"
$ws.Range("U210").Value = "g/m3"
$ws.Range("A211").Value = "Waitangi at d/s Waiouru STP"
$ws.Range("B211").Value = "Total Phosphorus (Median)"
$ws.Range("D211").Value = "2019 - 2023"
$ws.Range("E211").Value = "Impact"
$ws.Range("F211").Value = 0.037
$ws.Range("G211").Value = 0.0489152542372881
$ws.Range("H211").Value = 0.317
$ws.Range("I211").Value = 0.09594999999999999
$ws.Range("L211").Value = 0.035
$ws.Range("M211").Value = 0.05782
$ws.Range("N211").Value = 0.07056
$ws.Range("O211").Value = 1828659.395
$ws.Range("P211").Value = 5628533.029
$ws.Range("Q211").Value = "Ruapehu District"
$ws.Range("R211").Value = "Whangaehu"
$ws.Range("S211").Value = "Upper Whangaehu"
$ws.Range("T211").Value = "Whau_1b"
$ws.Range("U211").Value = "g/m3"

